$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $range = $ws.Range($Cell)
    $range.Value = "'" + $Text
    $range.Style = "Normal"
}

Set-TextValue "D2" "42.754.67"
Set-TextValue "E2" "  -1.98%  "
Set-TextValue "D3" "2.348.93"
Set-TextValue "E3" "  -2.88%  "
Set-TextValue "E4" "  -0.28%  "
Set-TextValue "D5" "320.68"
Set-TextValue "E5" "  -2.29%  "
Set-TextValue "D6" "105.16"
Set-TextValue "E6" "  +0.06%  "
Set-TextValue "E7" "  -1.96%  "
Set-TextValue "E8" "  -0.04%  "
Set-TextValue "E9" "  -8.20%  "
Set-TextValue "D10" "41.20"
Set-TextValue "E10" "  -2.96%  "
Set-TextValue "D11" "0.0922"
Set-TextValue "E11" "  -2.67%  "
Set-TextValue "D12" "8.43"
Set-TextValue "E12" "  -2.94%  "
Set-TextValue "E13" "  -3.16%  "
Set-TextValue "E14" "  -0.09%  "
Set-TextValue "D15" "15.97"
Set-TextValue "E15" "  -7.88%  "
Set-TextValue "D16" "2.702.80"
Set-TextValue "E16" "  -3.07%  "
Set-TextValue "D17" "2.362.80"
Set-TextValue "E17" "  -2.45%  "
Set-TextValue "D18" "42.718.73"
Set-TextValue "E18" "  -2.19%  "
Set-TextValue "E19" "  +1.60%  "
Set-TextValue "E20" "  -3.65%  "
Set-TextValue "D21" "77.06"
Set-TextValue "E21" "  +1.13%  "
Set-TextValue "D22" "3.61"
Set-TextValue "E22" "  +2.25%  "
Set-TextValue "D23" "261.21"
Set-TextValue "E23" "  -4.74%  "
Set-TextValue "E24" "  -6.04%  "
Set-TextValue "D25" "9.56"
Set-TextValue "E25" "  -0.89%  "
Set-TextValue "E26" "  +0.09%  "
Set-TextValue "D27" "11.37"
Set-TextValue "E27" "  -5.80%  "
Set-TextValue "D28" "23.10"
Set-TextValue "E28" "  -0.13%  "
Set-TextValue "E29" "  -0.28%  "
Set-TextValue "D30" "174.81"
Set-TextValue "E30" "  -1.99%  "
Set-TextValue "D31" "36.19"
Set-TextValue "E31" "  -4.79%  "
Set-TextValue "B32" "Filecoin"
Set-TextValue "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D32" "6.15"
Set-TextValue "E32" "  +2.68%  "
Set-TextValue "B33" "WEMIXToken"
Set-TextValue "C33" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D33" "2.99"
Set-TextValue "E33" "  -7.15%  "
Set-TextValue "D34" "0.0887"
Set-TextValue "E34" "  -5.70%  "
Set-TextValue "E35" "  -2.82%  "
Set-TextValue "D36" "0.118"
Set-TextValue "E36" "  +7.69%  "
Set-TextValue "D37" "4.62"
Set-TextValue "E37" "  -5.90%  "
Set-TextValue "E38" "  -3.61%  "
Set-TextValue "D39" "3.80"
Set-TextValue "E39" "  -6.66%  "
Set-TextValue "D40" "2.68"
Set-TextValue "E40" "  -7.60%  "
Set-TextValue "D41" "71.92"
Set-TextValue "E41" "  +2.22%  "
Set-TextValue "E42" "  -10.51%  "
Set-TextValue "E43" "  -2.26%  "
Set-TextValue "E44" "  -0.23%  "
Set-TextValue "D45" "115.33"
Set-TextValue "E45" "  -10.76%  "
Set-TextValue "D46" "89.14"
Set-TextValue "E46" "  +2.02%  "
Set-TextValue "D47" "11.87"
Set-TextValue "E47" "  -7.36%  "
Set-TextValue "D48" "5.49"
Set-TextValue "E48" "  -4.17%  "
Set-TextValue "E49" "  -6.67%  "
Set-TextValue "D50" "73.16"
Set-TextValue "E50" "  -0.14%  "
Set-TextValue "E51" "  -5.30%  "
